# "duplicate container data merged"
#
# Sheet "Laden Stock Report" had three consecutive rows (6, 7, 8) that all
# referred to the same container (PCIU2129609 / EIR 2014-3581) - rows 7 and 8
# were duplicate shipment/container records that only differed in their
# damage descriptions. The edit merges the duplicate container metadata away:
#   - row 6 keeps all of its data; its Issue/In Date cells (X6/Y6) pick up
#     the date number-format already used by the other rows (style index 4).
#   - rows 7 and 8 drop all of the duplicated container/shipment columns,
#     keep only the damage-description text (which moves from column O into
#     column N), and leave the Issue/In Date cells blank but still carrying
#     the date format.
#   - every row below (9..69) has its running serial number (column A)
#     reduced by 2 to account for the two rows that no longer carry unique
#     entries.
#   - columns N and O are resized to fit the new contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Laden Stock Report")

# --- Row 6: give the Issue/In Date cells the date number format ------------
$ws.Range("X6:Y6").NumberFormat = "YYYY-MM-DD"

# --- Row 7: keep only the damage description (former column O), in N ------
$descriptionRow7 = $ws.Cells.Item(7, 15).Value2
$ws.Range("A7:AG7").ClearContents()
$ws.Cells.Item(7, 14).Value = $descriptionRow7
$ws.Range("X7:Y7").NumberFormat = "YYYY-MM-DD"

# --- Row 8: keep only the damage description (former column O), in N ------
$descriptionRow8 = $ws.Cells.Item(8, 15).Value2
$ws.Range("A8:AG8").ClearContents()
$ws.Cells.Item(8, 14).Value = $descriptionRow8
$ws.Range("X8:Y8").NumberFormat = "YYYY-MM-DD"

# --- Renumber the serial-number column for every following row ------------
for ($r = 9; $r -le 69; $r++) {
    $current = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $current - 2
}

# --- Resize columns N and O to fit the merged content ----------------------
$ws.Columns.Item(14).ColumnWidth = 43.142857142857146
$ws.Columns.Item(15).ColumnWidth = 34.285714285714285
